# Update the build timestamp embedded in the version string throughout the
# workbook, replacing "February 03 2026 17.29.55 EST" with
# "February 03 2026 18.05.36 EST".

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: version banner (A2) and recommended citation (A6).
$cell = $wsAbout.Range("A2")
$cell.Value2 = $cell.Value2.Replace($oldTimestamp, $newTimestamp)

$cell = $wsAbout.Range("A6")
$cell.Value2 = $cell.Value2.Replace($oldTimestamp, $newTimestamp)

# Data sheet: build_version column (S) for every data row (2-33).
$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S
    if ($cell.Value2 -ne $null -and $cell.Value2 -is [string] -and $cell.Value2.Contains($oldTimestamp)) {
        $cell.Value2 = $cell.Value2.Replace($oldTimestamp, $newTimestamp)
    }
}
